$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update creature id column (A5:A16) from 2001xxx to 10xxx
$ws.Range("A5").Value = 10001
$ws.Range("A6").Value = 10002
$ws.Range("A7").Value = 10003
$ws.Range("A8").Value = 10004
$ws.Range("A9").Value = 10005
$ws.Range("A10").Value = 10006
$ws.Range("A11").Value = 10007
$ws.Range("A12").Value = 10008
$ws.Range("A13").Value = 10009
$ws.Range("A14").Value = 10010
$ws.Range("A15").Value = 10011
$ws.Range("A16").Value = 10012

# Update team/group column D (D6:D16) to 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1

# Update the selection to match the new active selection state
$ws.Range("A5:A16").Select()
